$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "datos actualizados" timestamp in A1
$ws.Range("A1").Value = "Datos actualizados a 25 de Marzo de 2020 a las 17:46"

# Update country rows: new case data causes a few countries to leapfrog
# their neighbours in the (descending, by total cases) ranking, and other
# rows get simple count updates.

# Row 39: Grecia
$ws.Range("A39").Value = "Grecia"
$ws.Range("B39").Value = 821
$ws.Range("C39").Value = 78
$ws.Range("D39").Value = 36
$ws.Range("E39").Value = 763
$ws.Range("F39").Value = 53
$ws.Range("G39").Value = 2
$ws.Range("H39").Value = 22

# Row 40: Indonesia
$ws.Range("A40").Value = "Indonesia"
$ws.Range("B40").Value = 790
$ws.Range("C40").Value = 104
$ws.Range("D40").Value = 31
$ws.Range("E40").Value = 701
$ws.Range("F40").Value = 0
$ws.Range("G40").Value = 3
$ws.Range("H40").Value = 58

# Row 51: Croacia
$ws.Range("A51").Value = "Croacia"
$ws.Range("B51").Value = 442
$ws.Range("C51").Value = 60
$ws.Range("D51").Value = 22
$ws.Range("E51").Value = 419
$ws.Range("F51").Value = 6
$ws.Range("G51").Value = 0
$ws.Range("H51").Value = 1

# Row 52: Egipto
$ws.Range("A52").Value = "Egipto"
$ws.Range("B52").Value = 442
$ws.Range("C52").Value = 40
$ws.Range("D52").Value = 93
$ws.Range("E52").Value = 328
$ws.Range("F52").Value = 0
$ws.Range("G52").Value = 1
$ws.Range("H52").Value = 21

# Row 53: Barein
$ws.Range("A53").Value = "Barein"
$ws.Range("B53").Value = 419
$ws.Range("C53").Value = 27
$ws.Range("D53").Value = 177
$ws.Range("E53").Value = 239
$ws.Range("F53").Value = 2
$ws.Range("G53").Value = 0
$ws.Range("H53").Value = 3

# Row 87: Republica de Chipre
$ws.Range("A87").Value = "Republica de Chipre"
$ws.Range("B87").Value = 132
$ws.Range("C87").Value = 8
$ws.Range("D87").Value = 3
$ws.Range("E87").Value = 126
$ws.Range("F87").Value = 3
$ws.Range("G87").Value = 0
$ws.Range("H87").Value = 3

# Row 88: Islas Feroe
$ws.Range("A88").Value = "Islas Feroe"
$ws.Range("B88").Value = 132
$ws.Range("C88").Value = 10
$ws.Range("D88").Value = 38
$ws.Range("E88").Value = 94
$ws.Range("F88").Value = 2
$ws.Range("G88").Value = 0
$ws.Range("H88").Value = 0

# Row 89: Malta
$ws.Range("A89").Value = "Malta"
$ws.Range("B89").Value = 129
$ws.Range("C89").Value = 19
$ws.Range("D89").Value = 2
$ws.Range("E89").Value = 127
$ws.Range("F89").Value = 1
$ws.Range("G89").Value = 0
$ws.Range("H89").Value = 0

# Row 6: Estados Unidos
$ws.Range("A6").Value = "Estados Unidos"
$ws.Range("B6").Value = 60567
$ws.Range("C6").Value = 5711
$ws.Range("D6").Value = 379
$ws.Range("E6").Value = 59372
$ws.Range("F6").Value = 1307
$ws.Range("G6").Value = 36
$ws.Range("H6").Value = 816

# Row 8: Alemania
$ws.Range("A8").Value = "Alemania"
$ws.Range("B8").Value = 35740
$ws.Range("C8").Value = 2749
$ws.Range("D8").Value = 3540
$ws.Range("E8").Value = 32014
$ws.Range("F8").Value = 23
$ws.Range("G8").Value = 27
$ws.Range("H8").Value = 186

# Row 65: Argelia
$ws.Range("A65").Value = "Argelia"
$ws.Range("B65").Value = 302
$ws.Range("C65").Value = 38
$ws.Range("D65").Value = 65
$ws.Range("E65").Value = 216
$ws.Range("F65").Value = 0
$ws.Range("G65").Value = 2
$ws.Range("H65").Value = 21

# Row 108: Estado de Palestina
$ws.Range("A108").Value = "Estado de Palestina"
$ws.Range("B108").Value = 64
$ws.Range("C108").Value = 4
$ws.Range("D108").Value = 16
$ws.Range("E108").Value = 48
$ws.Range("F108").Value = 0
$ws.Range("G108").Value = 0
$ws.Range("H108").Value = 0
